$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = "AllOntime"
$ws.Range("B3").Value = "PREC01049379"
$ws.Range("D1").Value = "Status"
$ws.Range("D3").Value = "PRECON"
$ws.Range("D1").Select()
